$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = New-Object 'object[,]' 4,34
$values[0,0] = 45123.50694444445
$values[0,1] = 0.596
$values[0,2] = 0.5659999999999999
$values[0,3] = 0.131
$values[0,4] = 0.703
$values[0,5] = 0.594
$values[0,6] = 0
$values[0,7] = 0.673
$values[0,8] = 2.927
$values[0,9] = 1.311
$values[0,10] = 0.574
$values[0,11] = 0.838
$values[0,12] = 0.092
$values[0,13] = 0.061
$values[0,14] = 0.781
$values[0,15] = 0.09
$values[0,16] = 0.143
$values[0,17] = 2.22
$values[0,18] = 0.626
$values[0,19] = 2.165
$values[0,20] = 0.66
$values[0,21] = 0.638
$values[0,22] = 1.735
$values[0,23] = 2.278
$values[0,24] = 0.128
$values[0,25] = 0.419
$values[0,26] = 0.467
$values[0,27] = 0.214
$values[0,28] = 0.476
$values[0,29] = 0.5600000000000001
$values[0,30] = 0.8090000000000001
$values[0,31] = 2.624
$values[0,32] = 0.591
$values[0,33] = 0.381
$values[1,0] = 45123.51388888889
$values[1,1] = 7.524
$values[1,2] = 5.666
$values[1,3] = 0.295
$values[1,4] = 16.004
$values[1,5] = 13.415
$values[1,6] = 5.44
$values[1,7] = 16.916
$values[1,8] = 9.797000000000001
$values[1,9] = 4.756
$values[1,10] = 6.212
$values[1,11] = 6.703
$values[1,12] = 6.736
$values[1,13] = 1.885
$values[1,14] = 6.166
$values[1,15] = 8.012
$values[1,16] = 4.56
$values[1,17] = 1.21
$values[1,18] = 0.347
$values[1,19] = 83.675
$values[1,20] = 15.969
$values[1,21] = 5.684
$values[1,22] = 11.477
$values[1,23] = 6.495
$values[1,24] = 0.793
$values[1,25] = 8.952
$values[1,26] = 4.887
$values[1,27] = 4.184
$values[1,28] = 5.047
$values[1,29] = 6.985
$values[1,30] = 0.343
$values[1,31] = 15.62
$values[1,32] = 3.401
$values[1,33] = 6.665
$values[2,0] = 45123.52083333334
$values[2,1] = 17.17
$values[2,2] = 12.908
$values[2,3] = 0.599
$values[2,4] = 37.146
$values[2,5] = 30.856
$values[2,6] = 13.187
$values[2,7] = 48.749
$values[2,8] = 21.224
$values[2,9] = 9.835000000000001
$values[2,10] = 14.042
$values[2,11] = 15.094
$values[2,12] = 15.739
$values[2,13] = 4.328
$values[2,14] = 13.664
$values[2,15] = 18.986
$values[2,16] = 10.885
$values[2,17] = 0.9330000000000001
$values[2,18] = 0.535
$values[2,19] = 198.398
$values[2,20] = 37.393
$values[2,21] = 12.615
$values[2,22] = 25.806
$values[2,23] = 13.766
$values[2,24] = 1.776
$values[2,25] = 24.653
$values[2,26] = 11.029
$values[2,27] = 9.617000000000001
$values[2,28] = 11.409
$values[2,29] = 15.844
$values[2,30] = 0.223
$values[2,31] = 44.83
$values[2,32] = 7.304
$values[2,33] = 15.424
$values[3,0] = 45123.52777777778
$values[3,1] = 22.48
$values[3,2] = 16.89
$values[3,3] = 0.77
$values[3,4] = 48.76
$values[3,5] = 40.43
$values[3,6] = 17.45
$values[3,7] = 68.09999999999999
$values[3,8] = 27.54
$values[3,9] = 12.62
$values[3,10] = 18.34
$values[3,11] = 19.71
$values[3,12] = 20.68
$values[3,13] = 5.67
$values[3,14] = 17.78
$values[3,15] = 25
$values[3,16] = 14.37
$values[3,17] = 0.79
$values[3,18] = 0.66
$values[3,19] = 261.45
$values[3,20] = 49.12
$values[3,21] = 16.42
$values[3,22] = 33.64
$values[3,23] = 17.79
$values[3,24] = 2.32
$values[3,25] = 33.51
$values[3,26] = 14.4
$values[3,27] = 12.6
$values[3,28] = 14.9
$values[3,29] = 20.71
$values[3,30] = 0.17
$values[3,31] = 62.29
$values[3,32] = 9.44
$values[3,33] = 20.24

$ws.Range("A2:AH5").Value = $values

$ws.Columns.Item(3).ColumnWidth = 7.17
$ws.Columns.Item(5).ColumnWidth = 7.17
$ws.Columns.Item(6).ColumnWidth = 7.17
$ws.Columns.Item(7).ColumnWidth = 7.17
$ws.Columns.Item(8).ColumnWidth = 7.17
$ws.Columns.Item(9).ColumnWidth = 7.17
$ws.Columns.Item(11).ColumnWidth = 7.17
$ws.Columns.Item(12).ColumnWidth = 7.17
$ws.Columns.Item(13).ColumnWidth = 7.17
$ws.Columns.Item(15).ColumnWidth = 7.17
$ws.Columns.Item(16).ColumnWidth = 7.17
$ws.Columns.Item(17).ColumnWidth = 7.17
$ws.Columns.Item(20).ColumnWidth = 8.17
$ws.Columns.Item(21).ColumnWidth = 7.17
$ws.Columns.Item(22).ColumnWidth = 7.17
$ws.Columns.Item(23).ColumnWidth = 7.17
$ws.Columns.Item(24).ColumnWidth = 7.17
$ws.Columns.Item(26).ColumnWidth = 7.17
$ws.Columns.Item(27).ColumnWidth = 7.17
$ws.Columns.Item(29).ColumnWidth = 7.17
$ws.Columns.Item(30).ColumnWidth = 7.17
$ws.Columns.Item(34).ColumnWidth = 7.17

$ws.Rows.Item(6).Delete()